# Insert a new data row before the current row 71 so every existing row
# from 71..152 shifts down by one (to 72..153), then populate the newly
# inserted row 71 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71 - this pushes rows 71..152 down to 72..153
# and extends the used range to A1:R153.
$ws.Rows.Item(71).Insert()

# Fill in the new row 71 with the new record's values (same shape/style as
# the other data rows in this table).
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44587
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112037
$ws.Cells.Item(71, 7).Value = "Cebollín"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 2800
$ws.Cells.Item(71, 11).Value = 900
$ws.Cells.Item(71, 12).Value = 1000
$ws.Cells.Item(71, 13).Value = 950
$ws.Cells.Item(71, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(71, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(71, 16).Value = 158
$ws.Cells.Item(71, 17).Value = 6
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same number format style ("s=2" in the
# original file) as the other dates in column D.
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
